$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.649.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "'3.432.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'407.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'130.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +2.52%  "
$ws.Range("D10").Value = "'0.139"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.89%  "
$ws.Range("D11").Value = "'42.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "'8.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "'19.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'3.444.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'62.544.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "'11.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "'1.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").Value = "'0.0000157"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +16.62%  "
$ws.Range("D20").Value = "'3.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").Value = "'84.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").Value = "'12.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("D24").Value = "'3.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").Value = "'29.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "'8.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.79%  "
$ws.Range("E28").Value = "  +4.05%  "
$ws.Range("E29").Value = "  +4.56%  "
$ws.Range("D30").Value = "'44.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "'0.114"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "'0.0483"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("D36").Value = "'51.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("E39").Value = "  +12.85%  "
$ws.Range("E40").Value = "  -4.21%  "
$ws.Range("D41").Value = "'142.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.26%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").Value = "'3.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").Value = "'16.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").Value = "'21.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("D48").Value = "'2.103.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").Value = "'1.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").Value = "'1.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +29.02%  "
